# Updates the cryptos list sheet: refreshed Price / Volume(1h) figures,
# plus two row swaps (dogwifhat <-> NEARProtocol, Maker <-> Fetch.AI)
# that happened because the underlying ranking reshuffled.
#
# Price values (column D) are stored as literal text in the workbook (e.g.
# "604.72", "69.967.43") rather than numbers, so every Price write below is
# prefixed with a leading apostrophe -- exactly what typing the value into
# Excel by hand would do -- to keep Excel from reinterpreting it as a
# number (which would silently drop things like trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($row, $text) {
    $ws.Cells.Item($row, 4).Value = "'" + $text
}

# --- Rows whose Price (D) and Volume(1h) (E) both change -------------------
$priceVolUpdates = @(
    @{ Row = 2;  D = "70.028.46";  E = "  +0.00%  " },
    @{ Row = 3;  D = "3.568.37";   E = "  +1.38%  " },
    @{ Row = 5;  D = "604.72";     E = "  -0.31%  " },
    @{ Row = 6;  D = "196.08";     E = "  -1.12%  " },
    @{ Row = 7;  D = "0.624";      E = "  -0.56%  " },
    @{ Row = 10; D = "0.650";      E = "  -1.48%  " },
    @{ Row = 11; D = "53.79";      E = "  -0.81%  " },
    @{ Row = 13; D = "9.55";       E = "  -1.21%  " },
    @{ Row = 14; D = "4.117.11";   E = "  +1.05%  " },
    @{ Row = 15; D = "597.23";     E = "  -0.69%  " },
    @{ Row = 16; D = "12.95";      E = "  +1.90%  " },
    @{ Row = 17; D = "19.27";      E = "  +1.09%  " },
    @{ Row = 18; D = "70.148.61";  E = "  +0.00%  " },
    @{ Row = 19; D = "3.558.07";   E = "  +0.79%  " },
    @{ Row = 20; D = "0.123";      E = "  +1.82%  " },
    @{ Row = 21; D = "0.991";      E = "  -1.07%  " },
    @{ Row = 22; D = "17.85";      E = "  -0.20%  " },
    @{ Row = 25; D = "4.64";       E = "  -0.53%  " },
    @{ Row = 26; D = "3.07";       E = "  -1.44%  " },
    @{ Row = 27; D = "10.82";      E = "  -1.82%  " },
    @{ Row = 28; D = "9.57";       E = "  -3.09%  " },
    @{ Row = 29; D = "33.57";      E = "  -1.25%  " },
    @{ Row = 32; D = "12.35";      E = "  -3.10%  " },
    @{ Row = 34; D = "63.44";      E = "  -0.84%  " },
    @{ Row = 37; D = "0.0₃0820";   E = "  +2.29%  " },
    @{ Row = 39; D = "519.73";     E = "  -0.50%  " },
    @{ Row = 40; D = "0.393";      E = "  +0.04%  " },
    @{ Row = 41; D = "3.61";       E = "  +1.09%  " },
    @{ Row = 42; D = "36.76";      E = "  -0.71%  " },
    @{ Row = 44; D = "0.0452";     E = "  -2.54%  " },
    @{ Row = 45; D = "2.84";       E = "  -1.35%  " },
    @{ Row = 47; D = "3.31";       E = "  +0.05%  " },
    @{ Row = 48; D = "8.55";       E = "  -3.05%  " },
    @{ Row = 50; D = "0.000247";   E = "  +2.81%  " }
)

foreach ($u in $priceVolUpdates) {
    Set-PriceText $u.Row $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# --- Rows whose Volume(1h) (E) alone changes --------------------------------
$volOnlyUpdates = @(
    @{ Row = 4;  E = "  -0.02%  " },
    @{ Row = 9;  E = "  -4.44%  " },
    @{ Row = 12; E = "  -0.70%  " },
    @{ Row = 23; E = "  +1.37%  " },
    @{ Row = 24; E = "  -1.33%  " },
    @{ Row = 33; E = "  -0.50%  " },
    @{ Row = 38; E = "  +0.17%  " },
    @{ Row = 43; E = "  -2.51%  " },
    @{ Row = 46; E = "  -0.18%  " },
    @{ Row = 49; E = "  +0.00%  " },
    @{ Row = 51; E = "  +3.09%  " }
)

foreach ($u in $volOnlyUpdates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# --- Row 30 / 31 swap: dogwifhat <-> NEARProtocol ---------------------------
$ws.Cells.Item(30, 2).Value = "NEARProtocol"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-PriceText 30 "7.11"
$ws.Cells.Item(30, 5).Value = "  -1.44%  "

$ws.Cells.Item(31, 2).Value = "dogwifhat"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-PriceText 31 "4.33"
$ws.Cells.Item(31, 5).Value = "  -3.93%  "

# --- Row 35 / 36 swap: Maker <-> Fetch.AI -----------------------------------
$ws.Cells.Item(35, 2).Value = "Fetch.AI"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-PriceText 35 "3.24"
$ws.Cells.Item(35, 5).Value = "  +7.36%  "

$ws.Cells.Item(36, 2).Value = "Maker"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-PriceText 36 "3.826.78"
$ws.Cells.Item(36, 5).Value = "  +2.83%  "
